{"js": "// Problem 1 - Updated solution and theory\n//\n// This script performs two independent OOXML-level paragraph replacements\n// that correspond to the commit's unified diff:\n//\n//  1) The paragraph beginning \"The only potential solutions are to leave\n//     the cat and the seed together...\" gets its final sentence appended\n//     in place, the trailing \"_GoBack\" bookmark + its old closing run are\n//     removed, and a large new block of analysis/discussion paragraphs is\n//     inserted, ending with a new paragraph that now carries the\n//     \"_GoBack\" bookmark.\n//\n//  2) The \"b) What insight...\" / \"the word problem alone?\" paragraphs get\n//     proofErr (grammar-check) markers added around \"from\" / \"the\", which\n//     also splits a couple of runs.\n\nfunction wrapFlatOpc(bodyInnerXml) {\n  return (\n    '<?xml version=\"1.0\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" + bodyInnerXml + \"<w:sectPr/></w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\n// ---------------------------------------------------------------------\n// Edit 1: replace the \"only potential solution(s)\" paragraph with the\n// extended sentence plus all of the new follow-up paragraphs.\n// ---------------------------------------------------------------------\nconst anchor1 = context.document.body.search(\n  \"the cat and the seed together without one of the items being consumed\",\n  { matchCase: true }\n);\nanchor1.load(\"text\");\nawait context.sync();\n\nif (anchor1.items.length === 0) {\n  throw new Error(\"Could not find the 'cat and seed' paragraph to update.\");\n}\n\nconst para1 = anchor1.items[0].paragraphs.getFirst();\nconst para1Range = para1.getRange(\"Whole\");\n\nconst newPara1BodyXml =\n  \"<w:p><w:r><w:t>The only</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> potential</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> solution</w:t></w:r>' +\n  \"<w:r><w:t>s are</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> to leave</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> the cat and the seed together without one of the items being consumed.  ' +\n  \"There are not any other options to pair the items together.</w:t></w:r></w:p>\" +\n  \"<w:p/>\" +\n  \"<w:p>\" +\n  \"<w:r><w:t>The potential and only solution would not meet the goals for the man.  \" +\n  \"Once the man took the parrot across the river he would leave the cat and seed behind.  \" +\n  \"His next trip would involve taking the cat or seed with him.  \" +\n  \"Which in either event would end with an item being consumed.  If he</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\"> took the seed on the second trip, it would be consumed by the parrot when he went back to get the cat.  ' +\n  \"If he took the cat on the second trip, the cat would consume the parrot when going back </w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\">for the seed.  The only option </w:t></w:r>' +\n  \"<w:r><w:t>the man has is to lose one of the items.</w:t></w:r>\" +\n  \"</w:p>\" +\n  \"<w:p/>\" +\n  \"<w:p><w:r><w:t>In order to test the accuracy of the theory above the break down of each avenue is below:</w:t></w:r></w:p>\" +\n  \"<w:p/>\" +\n  \"<w:p><w:r><w:t>Cat = c;</w:t></w:r></w:p>\" +\n  \"<w:p><w:r><w:t>Parrot = p;</w:t></w:r></w:p>\" +\n  \"<w:p><w:r><w:t>Seed = s;</w:t></w:r></w:p>\" +\n  \"<w:p/>\" +\n  \"<w:p><w:r><w:t>In order for the solution to work c &amp; p or p &amp; s can&#8217;t be together.</w:t></w:r></w:p>\" +\n  \"<w:p/>\" +\n  \"<w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\">Man takes c leaving </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  \"<w:r><w:t>p&amp;s</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> together = Unsuccessful</w:t></w:r>' +\n  \"</w:p>\" +\n  \"<w:p/>\" +\n  \"<w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\">Man takes p leaving </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  \"<w:r><w:t>c&amp;s</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> together = Successful</w:t></w:r>' +\n  \"</w:p>\" +\n  \"<w:p><w:r><w:t>Man takes c leaving s behind and going to p = Successful</w:t></w:r></w:p>\" +\n  \"<w:p><w:r><w:t>Man leaves c &amp; p together and going to s = Unsuccessful</w:t></w:r></w:p>\" +\n  \"<w:p/>\" +\n  \"<w:p><w:r><w:t>Man take s leaving c &amp; p together = Unsuccessful</w:t></w:r></w:p>\" +\n  \"<w:p/>\" +\n  \"<w:p>\" +\n  \"<w:r><w:lastRenderedPageBreak/><w:t>The 2</w:t></w:r>\" +\n  '<w:r><w:rPr><w:vertAlign w:val=\"superscript\"/></w:rPr><w:t>nd</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> theory has a shot at being successful initially, but is not able to be completed as there is no way to get through the three items without leaving two together that would consume one.</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  \"</w:p>\";\n\npara1Range.insertOoxml(wrapFlatOpc(newPara1BodyXml), Word.InsertLocation.replace);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Edit 2: add proofErr (grammar) markers around \"from\" in the \"b) What\n// insight...\" paragraph, and around \"the\" in the \"the word problem\n// alone?\" paragraph.\n// ---------------------------------------------------------------------\nconst anchor2a = context.document.body.search(\n  \"b) What insight can you offer into the problem that is not immediately\",\n  { matchCase: true }\n);\nanchor2a.load(\"text\");\nconst anchor2b = context.document.body.search(\"the word problem\", {\n  matchCase: true,\n});\nanchor2b.load(\"text\");\nawait context.sync();\n\nif (anchor2a.items.length === 0 || anchor2b.items.length === 0) {\n  throw new Error(\"Could not find the 'b) What insight' / 'the word problem' paragraphs to update.\");\n}\n\nconst paraB = anchor2a.items[0].paragraphs.getFirst();\nconst paraWordProblem = anchor2b.items[0].paragraphs.getFirst();\nconst twoParaRange = paraB.getRange(\"Start\").expandTo(paraWordProblem.getRange(\"End\"));\n\nconst newTwoParaBodyXml =\n  \"<w:p>\" +\n  \"<w:r><w:t>b) What insight can you offer into the problem that is not immediately</w:t></w:r>\" +\n  \"<w:r><w:tab/><w:t>visible</w:t></w:r>\" +\n  \"<w:r><w:tab/></w:r>\" +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  \"<w:r><w:t>from</w:t></w:r>\" +\n  \"<w:r><w:tab/></w:r>\" +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  \"</w:p>\" +\n  \"<w:p/>\" +\n  \"<w:p>\" +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  \"<w:r><w:t>the</w:t></w:r>\" +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> word problem </w:t></w:r>' +\n  \"<w:r><w:t>alone?</w:t></w:r>\" +\n  \"<w:r><w:tab/></w:r>\" +\n  \"</w:p>\";\n\ntwoParaRange.insertOoxml(wrapFlatOpc(newTwoParaBodyXml), Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Problem 1 - Updated solution and theory\n#\n# This script performs two independent OOXML-level paragraph replacements\n# that correspond to the commit's unified diff:\n#\n#  1) The paragraph beginning \"The only potential solutions are to leave\n#     the cat and the seed together...\" gets its final sentence appended\n#     in place, the trailing \"_GoBack\" bookmark + its old closing run are\n#     removed, and a large new block of analysis/discussion paragraphs is\n#     inserted, ending with a new paragraph that now carries the\n#     \"_GoBack\" bookmark.\n#\n#  2) The \"b) What insight...\" / \"the word problem alone?\" paragraphs get\n#     proofErr (grammar-check) markers added around \"from\" / \"the\", which\n#     also splits a couple of runs.\n\n$d = $word.ActiveDocument\n\nfunction Wrap-FlatOpc([string]$bodyInnerXml) {\n    return '<?xml version=\"1.0\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" ' +\n        'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $bodyInnerXml + '<w:sectPr/></w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n}\n\n# ---------------------------------------------------------------------\n# Edit 1: replace the \"only potential solution(s)\" paragraph with the\n# extended sentence plus all of the new follow-up paragraphs.\n# ---------------------------------------------------------------------\n$find1 = $d.Content\n$find1.Find.ClearFormatting()\n$found1 = $find1.Find.Execute(\"the cat and the seed together without one of the items being consumed\")\nif (-not $found1) {\n    throw \"Could not find the 'cat and seed' paragraph to update.\"\n}\n$para1Range = $find1.Paragraphs(1).Range\n\n$newPara1BodyXml = (\n    '<w:p><w:r><w:t>The only</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> potential</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> solution</w:t></w:r>' +\n    '<w:r><w:t>s are</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> to leave</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> the cat and the seed together without one of the items being consumed.  ' +\n    'There are not any other options to pair the items together.</w:t></w:r></w:p>' +\n    '<w:p/>' +\n    '<w:p>' +\n    '<w:r><w:t>The potential and only solution would not meet the goals for the man.  ' +\n    'Once the man took the parrot across the river he would leave the cat and seed behind.  ' +\n    'His next trip would involve taking the cat or seed with him.  ' +\n    'Which in either event would end with an item being consumed.  If he</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> took the seed on the second trip, it would be consumed by the parrot when he went back to get the cat.  ' +\n    'If he took the cat on the second trip, the cat would consume the parrot when going back </w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">for the seed.  The only option </w:t></w:r>' +\n    '<w:r><w:t>the man has is to lose one of the items.</w:t></w:r>' +\n    '</w:p>' +\n    '<w:p/>' +\n    '<w:p><w:r><w:t>In order to test the accuracy of the theory above the break down of each avenue is below:</w:t></w:r></w:p>' +\n    '<w:p/>' +\n    '<w:p><w:r><w:t>Cat = c;</w:t></w:r></w:p>' +\n    '<w:p><w:r><w:t>Parrot = p;</w:t></w:r></w:p>' +\n    '<w:p><w:r><w:t>Seed = s;</w:t></w:r></w:p>' +\n    '<w:p/>' +\n    '<w:p><w:r><w:t>In order for the solution to work c &amp; p or p &amp; s can&#8217;t be together.</w:t></w:r></w:p>' +\n    '<w:p/>' +\n    '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">Man takes c leaving </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>p&amp;s</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> together = Unsuccessful</w:t></w:r>' +\n    '</w:p>' +\n    '<w:p/>' +\n    '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">Man takes p leaving </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>c&amp;s</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> together = Successful</w:t></w:r>' +\n    '</w:p>' +\n    '<w:p><w:r><w:t>Man takes c leaving s behind and going to p = Successful</w:t></w:r></w:p>' +\n    '<w:p><w:r><w:t>Man leaves c &amp; p together and going to s = Unsuccessful</w:t></w:r></w:p>' +\n    '<w:p/>' +\n    '<w:p><w:r><w:t>Man take s leaving c &amp; p together = Unsuccessful</w:t></w:r></w:p>' +\n    '<w:p/>' +\n    '<w:p>' +\n    '<w:r><w:lastRenderedPageBreak/><w:t>The 2</w:t></w:r>' +\n    '<w:r><w:rPr><w:vertAlign w:val=\"superscript\"/></w:rPr><w:t>nd</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> theory has a shot at being successful initially, but is not able to be completed as there is no way to get through the three items without leaving two together that would consume one.</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    '</w:p>'\n)\n\n$para1Range.InsertXML((Wrap-FlatOpc $newPara1BodyXml))\n\n# ---------------------------------------------------------------------\n# Edit 2: add proofErr (grammar) markers around \"from\" in the \"b) What\n# insight...\" paragraph, and around \"the\" in the \"the word problem\n# alone?\" paragraph.\n# ---------------------------------------------------------------------\n$find2a = $d.Content\n$find2a.Find.ClearFormatting()\n$found2a = $find2a.Find.Execute(\"b) What insight can you offer into the problem that is not immediately\")\nif (-not $found2a) {\n    throw \"Could not find the 'b) What insight' paragraph to update.\"\n}\n$paraB = $find2a.Paragraphs(1).Range\n\n$find2b = $d.Content\n$find2b.Find.ClearFormatting()\n$found2b = $find2b.Find.Execute(\"the word problem\")\nif (-not $found2b) {\n    throw \"Could not find the 'the word problem' paragraph to update.\"\n}\n$paraWordProblem = $find2b.Paragraphs(1).Range\n\n$twoParaRange = $d.Range($paraB.Start, $paraWordProblem.End)\n\n$newTwoParaBodyXml = (\n    '<w:p>' +\n    '<w:r><w:t>b) What insight can you offer into the problem that is not immediately</w:t></w:r>' +\n    '<w:r><w:tab/><w:t>visible</w:t></w:r>' +\n    '<w:r><w:tab/></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>from</w:t></w:r>' +\n    '<w:r><w:tab/></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '</w:p>' +\n    '<w:p/>' +\n    '<w:p>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>the</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> word problem </w:t></w:r>' +\n    '<w:r><w:t>alone?</w:t></w:r>' +\n    '<w:r><w:tab/></w:r>' +\n    '</w:p>'\n)\n\n$twoParaRange.InsertXML((Wrap-FlatOpc $newTwoParaBodyXml))\n\nWrite-Output \"Edits applied.\"\n"}
